$wb = $excel.ActiveWorkbook

# --- Rename sheets (task-order tab names carry a new timestamp suffix) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16512554975008092"
$wb.Worksheets.Item(2).Name = "NB_TO-16512555006004918"
$wb.Worksheets.Item(3).Name = "RS_TO-16512555006064906"
$wb.Worksheets.Item(4).Name = "TOL_TO-16512555006644528"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16512555007414505"

# --- Sheet 1: GNG ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16512554974648116.csv"
$ws1.Range("B3").Value = "GNG_stims-16512554974838448.csv"
$ws1.Range("B4").Value = "go_stims-1651255497484844.csv"
$ws1.Range("B5").Value = "GNG_stims-1651255497499847.csv"

# --- Sheet 2: NB ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16512555005854964.csv"
$ws2.Range("B3").Value = "ZB-match_4-1651255498023817.csv"
$ws2.Range("B4").Value = "TB-1651255500399499.csv"
$ws2.Range("B5").Value = "OB-16512554993379278.csv"
$ws2.Range("B6").Value = "ZB-match_0-16512554979448082.csv"
$ws2.Range("B7").Value = "TB-16512554996394515.csv"
$ws2.Range("B8").Value = "OB-16512554989399862.csv"
$ws2.Range("B9").Value = "ZB-match_6-1651255498274851.csv"
$ws2.Range("B10").Value = "OB-16512554990949633.csv"

# --- Sheet 4: TOL ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16512555006324558.csv"
$ws4.Range("B3").Value = "ZM_stims-16512555006084547.csv"
$ws4.Range("B4").Value = "MM_stims-16512555006484556.csv"
$ws4.Range("B5").Value = "ZM_stims-16512555006334584.csv"
$ws4.Range("B6").Value = "MM_stims-16512555006634934.csv"
$ws4.Range("B7").Value = "ZM_stims-1651255500649457.csv"

# --- Sheet 5: vSAT ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-1651255500711452.csv"
$ws5.Range("B3").Value = "vSAT_stims-16512555007264504.csv"
$ws5.Range("B4").Value = "SAT_stims-1651255500695455.csv"
$ws5.Range("B5").Value = "SAT_stims-16512555006694658.csv"
